$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.ClearFormats()
}

Set-TextValue "D2" "27.365.33"
Set-TextValue "E2" "  +3.35%  "

Set-TextValue "D3" "1.750.44"
Set-TextValue "E3" "  +1.95%  "

Set-TextValue "D4" "0.9983"
Set-TextValue "E4" "  +0.41%  "

Set-TextValue "D5" "242.16"
Set-TextValue "E5" "  +0.79%  "

Set-TextValue "D6" "0.9984"
Set-TextValue "E6" "  +0.35%  "

Set-TextValue "D7" "0.4817"
Set-TextValue "E7" "  -1.66%  "

Set-TextValue "D8" "0.2614"
Set-TextValue "E8" "  +0.87%  "

Set-TextValue "D9" "0.06172"
Set-TextValue "E9" "  -0.18%  "

Set-TextValue "D10" "1.734.32"
Set-TextValue "E10" "  +0.95%  "

Set-TextValue "D11" "16.12"
Set-TextValue "E11" "  +3.56%  "

Set-TextValue "D12" "0.06939"
Set-TextValue "E12" "  -0.06%  "

Set-TextValue "D13" "0.6055"
Set-TextValue "E13" "  +0.47%  "

Set-TextValue "D14" "4.473"
Set-TextValue "E14" "  +0.29%  "

Set-TextValue "D15" "77.25"
Set-TextValue "E15" "  +0.89%  "

Set-TextValue "D16" "0.9983"
Set-TextValue "E16" "  +0.35%  "

Set-TextValue "D17" "27.345.12"
Set-TextValue "E17" "  +3.92%  "

Set-TextValue "D18" "0.9981"
Set-TextValue "E18" "  +0.41%  "

Set-TextValue "D19" "0.000007087"
Set-TextValue "E19" "  -0.26%  "

Set-TextValue "D20" "11.47"
Set-TextValue "E20" "  +1.80%  "

Set-TextValue "D21" "1.966.01"
Set-TextValue "E21" "  +1.69%  "

Set-TextValue "D22" "4.448"
Set-TextValue "E22" "  +1.22%  "

Set-TextValue "D23" "8.450"
Set-TextValue "E23" "  +0.65%  "

Set-TextValue "D24" "5.119"
Set-TextValue "E24" "  +1.24%  "

Set-TextValue "D25" "142.44"
Set-TextValue "E25" "  +3.67%  "

Set-TextValue "D26" "15.25"
Set-TextValue "E26" "  +0.29%  "

Set-TextValue "D27" "1.841"
Set-TextValue "E27" "  +6.11%  "

Set-TextValue "B28" "Toncoin"
Set-TextValue "C28" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D28" "1.410"
Set-TextValue "E28" "  +1.13%  "

Set-TextValue "B29" "BitcoinCash"
Set-TextValue "C29" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D29" "107.54"
Set-TextValue "E29" "  +2.29%  "

Set-TextValue "D30" "3.961"

Set-TextValue "D31" "0.07981"
Set-TextValue "E31" "  +0.30%  "

Set-TextValue "D32" "3.676"
Set-TextValue "E32" "  +1.51%  "

Set-TextValue "D33" "0.04684"
Set-TextValue "E33" "  +4.89%  "

Set-TextValue "E34" "  +0.30%  "

Set-TextValue "D35" "1.016"
Set-TextValue "E35" "  +1.78%  "

Set-TextValue "D36" "0.6192"
Set-TextValue "E36" "  +0.63%  "

Set-TextValue "D37" "0.9262"
Set-TextValue "E37" "  -3.20%  "

Set-TextValue "D38" "2.554"
Set-TextValue "E38" "  +7.85%  "

Set-TextValue "D39" "2.017"
Set-TextValue "E39" "  +1.29%  "

Set-TextValue "D40" "0.9987"
Set-TextValue "E40" "  +0.43%  "

Set-TextValue "D41" "5.735"
Set-TextValue "E41" "  +6.02%  "

Set-TextValue "D42" "0.01496"
Set-TextValue "E42" "  +1.17%  "

Set-TextValue "D43" "99.65"
Set-TextValue "E43" "  +0.06%  "

Set-TextValue "D44" "0.3853"
Set-TextValue "E44" "  +1.10%  "

Set-TextValue "D45" "6.903"
Set-TextValue "E45" "  +0.68%  "

Set-TextValue "D46" "0.1154"
Set-TextValue "E46" "  +0.35%  "

Set-TextValue "D47" "0.05363"
Set-TextValue "E47" "  +0.28%  "

Set-TextValue "D48" "7.849"
Set-TextValue "E48" "  +2.03%  "

Set-TextValue "D49" "29.89"
Set-TextValue "E49" "  -1.56%  "

Set-TextValue "D50" "1.250"
Set-TextValue "E50" "  +3.61%  "

Set-TextValue "D51" "51.09"
Set-TextValue "E51" "  -0.21%  "
